# Commit: "some exa linux vms"
# The DMZ VM entry (row 2: VSL-TST-OES-001 / nut-dmz-04 / ...) is removed
# from the "vms" sheet. Every row below it shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the whole row 2 (dmz_vms / VSL-TST-OES-001 / nut-dmz-04 ...),
# shifting rows 3:41 up to become rows 2:40.
$ws.Rows("2:2").Delete()

# The lone "sqref=H.." data validation keeps a *literal text* formula
# ("INDIRECT(G41)") that Excel does not auto-shift on row delete since it
# is not a real cell reference. Match what real Excel produces: re-point
# it at the new last data row (40) explicitly.
$ws.Range("H40").Validation.Modify(3, 1, 1, "INDIRECT(G40)")

# Restore the active selection to the cell the author ended up on.
$ws.Activate()
$ws.Range("E15").Select()
